$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet tab / title to reflect new date
$ws.Name = "Through 2022-05-06"

# Update header label cell (I1) text "2022 (through 05-05)" -> "2022 (through 05-06)"
$ws.Range("I1").Value = "2022 (through 05-06)"

# Update data values for May (row 6) and Total (row 14) in column I
$ws.Range("I6").Value = 19
$ws.Range("I14").Value = 570
